$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New precise values for rows 2-5, columns B-X (columns 2-24)
$data = @{
    2 = @(0.00950020652622883, 0.000619578686493185, 0.980896323833127, 0.0170384138785626, 0.00216852540272615, 0.877426683188765, 0.00103263114415531, 0.000413052457662123, 0.0219950433705081, 0.0026848409748038, 0.998864105741429, 0.000103263114415531, 0.000516315572077654, 0.000103263114415531, 0.969021065675341, 0.0026848409748038, 0.998657579512598, 0.961482858323007, 0.980380008261049, 0.00196199917389508, 0.993700950020653, 0.00454357703428335, 0.135997521685254)
    3 = @(0.000206526228831062, 0.000826104915324246, 0.000929368029739777, 0.000309789343246592, 0.0193102023957043, 0.0469847170590665, 0.000516315572077654, 0.998244527054936, 0.00309789343246592, 0.000413052457662123, 0, 0.870095002065262, 0.0365551425030979, 0, 0.000309789343246592, 0.000826104915324246, 0.000103263114415531, 0.00103263114415531, 0.00908715406856671, 0.00278810408921933, 0.000929368029739777, 0.000413052457662123, 0.000619578686493185)
    4 = @(0.986162742668319, 0.000103263114415531, 0.0114622057001239, 0.980276745146634, 0.974597273853779, 0.0577240809582817, 0.000516315572077654, 0.000619578686493185, 0.00660883932259397, 0.995043370508054, 0.000206526228831062, 0.000309789343246592, 0.000103263114415531, 0, 0.0263320941759603, 0.995662949194548, 0.000722841800908715, 0.028397356464271, 0.000826104915324246, 0.000413052457662123, 0.00175547294506402, 0.993391160677406, 0.85212722015696)
    5 = @(0.000413052457662123, 0.99814126394052, 0.00123915737298637, 0.00123915737298637, 0, 0.00454357703428335, 0.997934737711689, 0.000103263114415531, 0.960140437835605, 0, 0.000826104915324246, 0.111524163568773, 0.956732755059893, 0.999793473771169, 0.000826104915324246, 0, 0.000103263114415531, 0.00123915737298637, 0.00722841800908715, 0.993184634448575, 0.00258157786038827, 0.000413052457662123, 0.000516315572077654)
}

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $colNum = $i + 2  # Column B = 2
        $ws.Cells.Item($rowNum, $colNum).Value = $values[$i]
    }
}
